$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'29.687.16"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "'2.095.65"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'343.16"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "'0.5162"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.4378"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'53.57"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("D11").Value = "'1.165"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  -5.12%  "
$ws.Range("D13").Value = "'2.080.05"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "'6.762"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "'8.150"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'102.27"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "'0.00001151"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'21.00"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'0.06669"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'6.196"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'29.756.49"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").Value = "'12.63"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'2.312.15"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'21.89"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "'161.88"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "'2.493"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").Value = "'133.18"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "'1.128"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("D32").Value = "'1.664"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'0.1050"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "'6.185"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").Value = "'3.961"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'6.285"
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("D37").Value = "'10.41"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "'0.02574"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'0.06691"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Value = "'0.6980"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "'12.43"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "'1.328"
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("D43").Value = "'0.2212"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("D44").Value = "'0.6797"
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("D45").Value = "'14.29"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "'2.316"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'0.00000000363"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'3.623"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'1.205"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").Value = "'1.215"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "'81.10"
$ws.Range("E51").Value = "  -3.70%  "
